$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the current content of rows 4-33 (values + key formulas) ---
# Columns used in data rows: A..R (1..18) for plain values, S,T,V,W,X,Y (19,20,22,23,24,25) for HYPERLINK formulas.
$snapshot = @{}
for ($r = 4; $r -le 33; $r++) {
    $rowData = @{}
    for ($c = 1; $c -le 18; $c++) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    # Hyperlink formula columns: S=19, T=20, V=22, W=23, X=24, Y=25 (U=21, Z=26 unused in data rows)
    foreach ($c in 19,20,22,23,24,25) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $rowData[$c] = $cell.Formula
        } else {
            $rowData[$c] = $null
        }
    }
    $snapshot[$r] = $rowData
}

# --- Step 2: row-content relocation map (destination row -> source row), per the authoritative edit ---
$rowMap = @{
    4 = 5
    5 = 6
    6 = 4
    7 = 7
    8 = 9
    9 = 13
    10 = 16
    11 = 17
    12 = 19
    13 = 25
    14 = 8
    15 = 10
    16 = 12
    17 = 28
    18 = 29
    19 = 18
    20 = 21
    21 = 22
    22 = 15
    23 = 20
    24 = 24
    25 = 26
    26 = 32
    27 = 33
    28 = 27
    29 = 30
    30 = 31
    31 = 11
    32 = 14
    33 = 23
}

# --- Step 3: write each destination row from its mapped source snapshot ---
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $data = $snapshot[$srcRow]
    for ($c = 1; $c -le 18; $c++) {
        $v = $data[$c]
        if ($null -eq $v -or $v -eq "") {
            $ws.Cells.Item($destRow, $c).Value2 = ""
        } else {
            $ws.Cells.Item($destRow, $c).Value2 = $v
        }
    }
    foreach ($c in 19,20,22,23,24,25) {
        $f = $data[$c]
        if ($null -ne $f) {
            $ws.Cells.Item($destRow, $c).Formula = $f
        } else {
            $ws.Cells.Item($destRow, $c).ClearContents() | Out-Null
        }
    }
}

# --- Step 4: update column C ("Förändrad") for every data row 2-33 to the new date serial 46081 ---
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46081
}

Write-Output "done"